$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030361113059984
$ws.Range("D2").Value = 1.040077525611573
$ws.Range("E2").Value = 1.0300745788704
$ws.Range("F2").Value = 1.048174272068829
$ws.Range("I2").Value = 1.033256158706183
$ws.Range("J2").Value = 1.035502783351589
$ws.Range("K2").Value = 1.042860897267964
$ws.Range("L2").Value = 1.032886643825647
$ws.Range("M2").Value = 1.050934853409742
$ws.Range("N2").Value = 1.015815989273407

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031254538043011
$ws.Range("D3").Value = 1.040916154152011
$ws.Range("E3").Value = 1.030832069581703
$ws.Range("F3").Value = 1.049096986422159
$ws.Range("I3").Value = 1.033364926305235
$ws.Range("J3").Value = 1.036037915975943
$ws.Range("K3").Value = 1.043509800910734
$ws.Range("L3").Value = 1.033452549832663
$ws.Range("M3").Value = 1.051669273385551
$ws.Range("N3").Value = 1.015995216245234

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031833183624985
$ws.Range("D4").Value = 1.041459624285112
$ws.Range("E4").Value = 1.031323062490476
$ws.Range("F4").Value = 1.049694949048365
$ws.Range("I4").Value = 1.033433990902831
$ws.Range("J4").Value = 1.036384092389006
$ws.Range("K4").Value = 1.043929860109491
$ws.Range("L4").Value = 1.033818912507516
$ws.Range("M4").Value = 1.052144766924754
$ws.Range("N4").Value = 1.016111104140032

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032076574150214
$ws.Range("D5").Value = 1.041688294152839
$ws.Range("E5").Value = 1.031529676777438
$ws.Range("F5").Value = 1.049946547468012
$ws.Range("I5").Value = 1.033462710181774
$ws.Range("J5").Value = 1.036529602211118
$ws.Range("K5").Value = 1.04410649359469
$ws.Range("L5").Value = 1.033972974106336
$ws.Range("M5").Value = 1.052344728712729
$ws.Range("N5").Value = 1.016159802878386

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032117447953373
$ws.Range("D6").Value = 1.041726700219063
$ws.Range("E6").Value = 1.031564379957537
$ws.Range("F6").Value = 1.049988804505075
$ws.Range("I6").Value = 1.033467513758182
$ws.Range("J6").Value = 1.036554032583411
$ws.Range("K6").Value = 1.044136153476011
$ws.Range("L6").Value = 1.033998844196492
$ws.Range("M6").Value = 1.052378306905171
$ws.Range("N6").Value = 1.016167978394086

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031836435320099
$ws.Range("D7").Value = 1.041462679018797
$ws.Range("E7").Value = 1.03132582249344
$ws.Range("F7").Value = 1.049698310076316
$ws.Range("I7").Value = 1.033434375891588
$ws.Range("J7").Value = 1.036386036789177
$ws.Range("K7").Value = 1.043932220136915
$ws.Range("L7").Value = 1.033820970919789
$ws.Range("M7").Value = 1.052147438573113
$ws.Range("N7").Value = 1.016111754936185

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.030662937633419
$ws.Range("D8").Value = 1.040360773005781
$ws.Range("E8").Value = 1.030330400313516
$ws.Range("F8").Value = 1.048485919371039
$ws.Range("I8").Value = 1.033293189214971
$ws.Range("J8").Value = 1.035683651906098
$ws.Range("K8").Value = 1.043080160224206
$ws.Range("L8").Value = 1.033077855642336
$ws.Range("M8").Value = 1.051182996320867
$ws.Range("N8").Value = 1.015876576985419

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028599278186554
$ws.Range("D9").Value = 1.038425436220835
$ws.Range("E9").Value = 1.028582887943676
$ws.Range("F9").Value = 1.046356534866891
$ws.Range("I9").Value = 1.033034352591554
$ws.Range("J9").Value = 1.034445320287627
$ws.Range("K9").Value = 1.041580122160905
$ws.Range("L9").Value = 1.031769852663005
$ws.Range("M9").Value = 1.049485694335723
$ws.Range("N9").Value = 1.015461539801434

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.027226402913728
$ws.Range("D10").Value = 1.03713958549379
$ws.Range("E10").Value = 1.027422378361227
$ws.Range("F10").Value = 1.044941752270808
$ws.Range("I10").Value = 1.032855074365956
$ws.Range("J10").Value = 1.033619404126182
$ws.Range("K10").Value = 1.040581117398528
$ws.Range("L10").Value = 1.030898906334413
$ws.Range("M10").Value = 1.048355704305471
$ws.Range("N10").Value = 1.015184454994978

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026632636324666
$ws.Range("D11").Value = 1.036583855492023
$ws.Range("E11").Value = 1.026920952047105
$ws.Range("F11").Value = 1.044330296648812
$ws.Range("I11").Value = 1.032775858229259
$ws.Range("J11").Value = 1.033261701815371
$ws.Range("K11").Value = 1.040148795946533
$ws.Range("L11").Value = 1.03052204272348
$ws.Range("M11").Value = 1.047866790658417
$ws.Range("N11").Value = 1.015064386348217

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.026412191042092
$ws.Range("D12").Value = 1.036377592269138
$ws.Range("E12").Value = 1.026734864166497
$ws.Range("F12").Value = 1.044103349922134
$ws.Range("I12").Value = 1.032746195894456
$ws.Range("J12").Value = 1.033128825093853
$ws.Range("K12").Value = 1.03998825202028
$ws.Range("L12").Value = 1.030382099374002
$ws.Range("M12").Value = 1.047685244820718
$ws.Range("N12").Value = 1.015019774552883

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026459472483215
$ws.Range("D13").Value = 1.036421829180474
$ws.Range("E13").Value = 1.02677477319314
$ws.Range("F13").Value = 1.044152022814438
$ws.Range("I13").Value = 1.032752569323894
$ws.Range("J13").Value = 1.033157328043045
$ws.Range("K13").Value = 1.040022687428339
$ws.Range("L13").Value = 1.030412115841667
$ws.Range("M13").Value = 1.047724184350538
$ws.Range("N13").Value = 1.015029344513678

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026614412077706
$ws.Range("D14").Value = 1.036566802443651
$ws.Range("E14").Value = 1.026905566606281
$ws.Range("F14").Value = 1.044311533579696
$ws.Range("I14").Value = 1.032773411181922
$ws.Range("J14").Value = 1.03325071838474
$ws.Range("K14").Value = 1.040135524521882
$ws.Range("L14").Value = 1.030510474127726
$ws.Range("M14").Value = 1.047851782821918
$ws.Range("N14").Value = 1.015060698985056

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026709889568671
$ws.Range("D15").Value = 1.036656146458471
$ws.Range("E15").Value = 1.026986174548558
$ws.Range("F15").Value = 1.044409836680815
$ws.Range("I15").Value = 1.032786221026855
$ws.Range("J15").Value = 1.033308257944665
$ws.Range("K15").Value = 1.040205052454251
$ws.Range("L15").Value = 1.030571081327386
$ws.Range("M15").Value = 1.047930408241141
$ws.Range("N15").Value = 1.015080015806103

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.027265824050644
$ws.Range("D16").Value = 1.037176489762867
$ws.Range("E16").Value = 1.027455679309055
$ws.Range("F16").Value = 1.044982357033454
$ws.Range("I16").Value = 1.032860298281985
$ws.Range("J16").Value = 1.03364314216383
$ws.Range("K16").Value = 1.040609814617817
$ws.Range("L16").Value = 1.030923923153829
$ws.Range("M16").Value = 1.048388160012842
$ws.Range("N16").Value = 1.015192421715332

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027614734696765
$ws.Range("D17").Value = 1.037503170120354
$ws.Range("E17").Value = 1.027750478066358
$ws.Range("F17").Value = 1.045341794205512
$ws.Range("I17").Value = 1.032906340372786
$ws.Range("J17").Value = 1.033853186865433
$ws.Range("K17").Value = 1.040863780217605
$ws.Range("L17").Value = 1.031145322494085
$ws.Range("M17").Value = 1.048675398226506
$ws.Range("N17").Value = 1.015262907423675

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027818315582615
$ws.Range("D18").Value = 1.037693818625657
$ws.Range("E18").Value = 1.027922533481901
$ws.Range("F18").Value = 1.045551559115836
$ws.Range("I18").Value = 1.032933042738338
$ws.Range("J18").Value = 1.033975695057898
$ws.Range("K18").Value = 1.04101193849357
$ws.Range("L18").Value = 1.031274486170341
$ws.Range("M18").Value = 1.048842975935429
$ws.Range("N18").Value = 1.015304011919298

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027887742764628
$ws.Range("D19").Value = 1.037758841985794
$ws.Range("E19").Value = 1.027981217572004
$ws.Range("F19").Value = 1.045623102355697
$ws.Range("I19").Value = 1.03294212155114
$ws.Range("J19").Value = 1.034017465908078
$ws.Range("K19").Value = 1.041062460741969
$ws.Range("L19").Value = 1.031318531880807
$ws.Range("M19").Value = 1.048900121753808
$ws.Range("N19").Value = 1.015318026013885

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027577292925486
$ws.Range("D20").Value = 1.037468109911812
$ws.Range("E20").Value = 1.02771883815094
$ws.Range("F20").Value = 1.045303218475528
$ws.Range("I20").Value = 1.032901416335012
$ws.Range("J20").Value = 1.033830651806587
$ws.Range("K20").Value = 1.040836529574947
$ws.Range("L20").Value = 1.031121565826107
$ws.Range("M20").Value = 1.048644576496843
$ws.Range("N20").Value = 1.015255345863297

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026568783290654
$ws.Range("D21").Value = 1.03652410701109
$ws.Range("E21").Value = 1.026867046660909
$ws.Range("F21").Value = 1.044264556787669
$ws.Range("I21").Value = 1.032767280337805
$ws.Range("J21").Value = 1.033223217528052
$ws.Range("K21").Value = 1.040102295708955
$ws.Range("L21").Value = 1.030481508934052
$ws.Range("M21").Value = 1.04781420662702
$ws.Range("N21").Value = 1.015051466228847

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025935307557867
$ws.Range("D22").Value = 1.035931499368786
$ws.Range("E22").Value = 1.026332442231203
$ws.Range("F22").Value = 1.043612522949426
$ws.Range("I22").Value = 1.03268156737789
$ws.Range("J22").Value = 1.032841241287883
$ws.Range("K22").Value = 1.039640882982592
$ws.Range("L22").Value = 1.030079315127372
$ws.Range("M22").Value = 1.047292459078538
$ws.Range("N22").Value = 1.014923204277243

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.026271066372287
$ws.Range("D23").Value = 1.036245563725966
$ws.Range("E23").Value = 1.026615755499224
$ws.Range("F23").Value = 1.043958081756815
$ws.Range("I23").Value = 1.03272713571083
$ws.Range("J23").Value = 1.033043739257355
$ws.Range("K23").Value = 1.03988546442849
$ws.Range("L23").Value = 1.030292502954017
$ws.Range("M23").Value = 1.047569014714773
$ws.Range("N23").Value = 1.014991205304288

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027594211042722
$ws.Range("D24").Value = 1.03748395179784
$ws.Range("E24").Value = 1.027733134544753
$ws.Range("F24").Value = 1.045320648843751
$ws.Range("I24").Value = 1.032903641769175
$ws.Range("J24").Value = 1.033840834452434
$ws.Range("K24").Value = 1.040848842892008
$ws.Range("L24").Value = 1.031132300364087
$ws.Range("M24").Value = 1.048658503397091
$ws.Range("N24").Value = 1.015258762633648

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.02913227862665
$ws.Range("D25").Value = 1.038925003293464
$ws.Range("E25").Value = 1.029033875860576
$ws.Range("F25").Value = 1.046906191759198
$ws.Range("I25").Value = 1.03310245515895
$ws.Range("J25").Value = 1.034765527460789
$ws.Range("K25").Value = 1.041967743858008
$ws.Range("L25").Value = 1.032107822085173
$ws.Range("M25").Value = 1.049924221680216
$ws.Range("N25").Value = 1.01556890808884
